$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''57.694.01'
$ws.Range("E2").Value = '  +1.49%  '
$ws.Range("D3").Value = '''3.015.94'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''510.76'
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").Value = '''140.31'
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '''0.438'
$ws.Range("E8").Value = '  +0.78%  '
$ws.Range("D9").Value = '''7.59'
$ws.Range("E9").Value = '  -0.19%  '
$ws.Range("D10").Value = '''0.111'
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("D11").Value = '''0.367'
$ws.Range("E11").Value = '  +3.26%  '
$ws.Range("D12").Value = '''3.528.71'
$ws.Range("E12").Value = '  +0.41%  '
$ws.Range("E13").Value = '  +0.72%  '
$ws.Range("D14").Value = '''26.67'
$ws.Range("E14").Value = '  +3.33%  '
$ws.Range("D15").Value = '''0.0000164'
$ws.Range("E15").Value = '  +5.24%  '
$ws.Range("D16").Value = '''57.701.30'
$ws.Range("E16").Value = '  +1.51%  '
$ws.Range("D17").Value = '''6.24'
$ws.Range("E17").Value = '  +5.74%  '
$ws.Range("D18").Value = '''3.015.75'
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("D19").Value = '''12.91'
$ws.Range("E19").Value = '  +3.36%  '
$ws.Range("E20").Value = '  +1.10%  '
$ws.Range("D21").Value = '''331.75'
$ws.Range("E21").Value = '  +0.36%  '
$ws.Range("D22").Value = '''0.999'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").Value = '''0.500'
$ws.Range("E23").Value = '  +3.14%  '
$ws.Range("D24").Value = '''64.68'
$ws.Range("E24").Value = '  +3.04%  '
$ws.Range("E25").Value = '  -2.26%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").Value = '''0.0₃0927'
$ws.Range("E27").Value = '  +0.94%  '
$ws.Range("D28").Value = '''6.84'
$ws.Range("E28").Value = '  +2.31%  '
$ws.Range("D29").Value = '''7.33'
$ws.Range("E29").Value = '  +1.00%  '
$ws.Range("E30").Value = '  +1.58%  '
$ws.Range("D31").Value = '''1.20'
$ws.Range("E31").Value = '  -6.05%  '
$ws.Range("D32").Value = '''20.72'
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("D33").Value = '''4.77'
$ws.Range("E33").Value = '  +4.29%  '
$ws.Range("D34").Value = '''155.12'
$ws.Range("E34").Value = '  -0.84%  '
$ws.Range("D35").Value = '''5.90'
$ws.Range("E35").Value = '  +3.87%  '
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").Value = '''24.59'
$ws.Range("E37").Value = '  +2.17%  '
$ws.Range("D38").Value = '''0.0679'
$ws.Range("E38").Value = '  -0.61%  '
$ws.Range("D39").Value = '''3.048.78'
$ws.Range("E39").Value = '  +0.42%  '
$ws.Range("D40").Value = '''37.62'
$ws.Range("E40").Value = '  +1.80%  '
$ws.Range("D41").Value = '''3.87'
$ws.Range("E41").Value = '  +5.85%  '
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").Value = '''0.652'
$ws.Range("E43").Value = '  +0.72%  '
$ws.Range("D44").Value = '''1.42'
$ws.Range("E44").Value = '  -0.50%  '
$ws.Range("D45").Value = '''2.237.04'
$ws.Range("E45").Value = '  -1.63%  '
$ws.Range("D46").Value = '''0.989'
$ws.Range("E46").Value = '  -1.46%  '
$ws.Range("D47").Value = '''6.03'
$ws.Range("E47").Value = '  +3.80%  '
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").Value = '''19.58'
$ws.Range("E49").Value = '  +0.90%  '
$ws.Range("E50").Value = '  -8.54%  '
$ws.Range("D51").Value = '''0.0896'
$ws.Range("E51").Value = '  +2.97%  '

Write-Host "Applied cryptos update"
